# Commit: Trade #83 closed at 2026-02-17 21:17:42 - unknown UNKNOWN +0.000%
#
# - Updates the Summary sheet's aggregate stats
# - Updates the Strategy Status row for MarketMaking
# - Marks existing trade #111 as CLOSED (early_exit) in both
#   "All Trades" and "MarketMaking" sheets
# - Appends a brand-new OPEN trade #144 to both "All Trades" and
#   "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.14   # Current Capital
$summary.Range("B4").Value = 0.93      # Total P&L $
$summary.Range("B6").Value = 111       # Total Trades
$summary.Range("B7").Value = 50        # Winning Trades
$summary.Range("B9").Value = 45.05     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: Strategy Status (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.14     # Capital
$status.Range("D5").Value = 78         # Trades
$status.Range("E5").Value = 0.82       # P&L $
$status.Range("F5").Value = 1.14       # P&L %
$status.Range("G5").Value = 46.15      # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: All Trades - close out existing trade #111 (row 112)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G112").Value = 0.84
$allTrades.Range("H112").Value = "CLOSED"
$allTrades.Range("I112").Value = 1.2048
$allTrades.Range("J112").Value = 0.01
$allTrades.Range("K112").Value = 101.14
$allTrades.Range("L112").Value = "early_exit"
$allTrades.Range("M112").Value = 0.13

# Append new trade #144 (row 145)
$allTrades.Range("A145").Value = 144
$allTrades.Range("B145").NumberFormat = "@"
$allTrades.Range("B145").Value = "2026-02-17"
$allTrades.Range("B145").NumberFormat = "General"
$allTrades.Range("C145").Value = "21:17:35"
$allTrades.Range("D145").Value = "MarketMaking"
$allTrades.Range("E145").Value = "DOWN"
$allTrades.Range("F145").Value = 0.83
$allTrades.Range("H145").Value = "OPEN"
$allTrades.Range("I145").Value = 0
$allTrades.Range("J145").Value = 0
$allTrades.Range("K145").Value = 101.1296151053151
$allTrades.Range("M145").Value = 0
$allTrades.Range("N145").Value = 0
$allTrades.Range("O145").Value = 0
$allTrades.Range("P145").Value = 0.6
$allTrades.Range("Q145").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet: MarketMaking - close out existing trade #111 (row 79)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G79").Value = 0.84
$mm.Range("H79").Value = "CLOSED"
$mm.Range("I79").Value = 1.2048
$mm.Range("J79").Value = 0.01
$mm.Range("K79").Value = 101.14
$mm.Range("P79").Value = "early_exit"
$mm.Range("Q79").Value = 0.13

# Append new trade #144 (row 112)
$mm.Range("A112").Value = 144
$mm.Range("B112").NumberFormat = "@"
$mm.Range("B112").Value = "2026-02-17"
$mm.Range("B112").NumberFormat = "General"
$mm.Range("C112").Value = "21:17:35"
$mm.Range("D112").Value = "MarketMaking"
$mm.Range("E112").Value = "DOWN"
$mm.Range("F112").Value = 0.83
$mm.Range("H112").Value = "OPEN"
$mm.Range("I112").Value = 0
$mm.Range("J112").Value = 0
$mm.Range("K112").Value = 101.1296151053151
$mm.Range("L112").Value = 0
$mm.Range("M112").Value = 0
$mm.Range("N112").Value = 0.6
$mm.Range("O112").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q112").Value = 0
